# Update the CS (ClientServer) interface row on Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 8: PSDContrl / provide / SetDoorCmd / SetDoorCmd / SetDoorCmdOperation / ClientServer / uint8
# (string order matters for the shared-string table, so write ClientServer
# before SetDoorCmdOperation to match the author's original entry order)
$ws.Range("B8").Value = "PSDContrl"
$ws.Range("C8").Value = "provide"
$ws.Range("D8").Value = "SetDoorCmd"
$ws.Range("E8").Value = "SetDoorCmd"
$ws.Range("G8").Value = "ClientServer"
$ws.Range("F8").Value = "SetDoorCmdOperation"
$ws.Range("H8").Value = "uint8"

# Match style of the row above it (B2:H7 data rows all share one style).
$ws.Range("B2:H2").Copy()
$ws.Range("B8:H8").PasteSpecial(-4122)  # xlPasteFormats

# Widen column F slightly to fit the new ElementName text.
$ws.Columns("F").ColumnWidth = 20.75

# Update the active selection to reflect where the author left off editing.
$ws.Range("E15").Select()

# Configure the page setup (paper size / orientation) used when the sheet
# was last saved.
$ws.PageSetup.PaperSize = 9   # xlPaperA4
$ws.PageSetup.Orientation = 1 # xlPortrait
